$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.02.2022 18:00"

# Row 7 (MOL Olomoucká) price refresh:
# B7 gets the new price, C7 becomes the previous price, D7 is the textual
# delta (with explicit sign), E7 is a literal timestamp string.
$ws.Range("B7").Value = 38.29
$ws.Range("C7").Value = 37.9

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "+0.39"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2022-02-21 18:00:26"

# Drop the number-format overrides so D7/E7 end up with the default style
# (matches the target: no explicit "s" attribute on either cell).
$ws.Range("D7:E7").ClearFormats()
